$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "E2" = 3
    "F2" = 1
    "G2" = 1.342908333333333
    "H2" = 4.028725
    "I2" = 0.2879023314891748
    "J2" = 0.2879023314891748
    "M2" = 4.901461666666667
    "N2" = 14.704385
    "O2" = 0.2124427850531459
    "P2" = 0.2124427850531459
    "Q2" = 6.582213717680555
    "R2" = 59.239923459125
    "S2" = 0.0611627731248543
    "T2" = 0.06116277312485431
    "E3" = 3
    "F3" = 1
    "G3" = 1.342908333333333
    "H3" = 4.028725
    "I3" = 0.2879023314891748
    "J3" = 0.2879023314891748
    "O3" = 0.1372144215401173
    "P3" = 0.1372144215401173
    "Q3" = 4.251378306394444
    "R3" = 38.26240475755
    "S3" = 0.03950435187533821
    "T3" = 0.03950435187533821
    "E4" = 3
    "F4" = 1
    "G4" = 1.342908333333333
    "H4" = 4.028725
    "I4" = 0.2879023314891748
    "J4" = 0.2879023314891748
    "M4" = 1.206743666666667
    "N4" = 3.620231
    "O4" = 0.05230357857032003
    "P4" = 0.05230357857032004
    "Q4" = 1.620546126163889
    "R4" = 14.584915135475
    "S4" = 0.01505832221562237
    "T4" = 0.01505832221562238
    "E5" = 3
    "F5" = 1
    "G5" = 1.342908333333333
    "H5" = 4.028725
    "I5" = 0.2879023314891748
    "J5" = 0.2879023314891748
    "M5" = 13.79790933333333
    "N5" = 41.393728
    "O5" = 0.5980392148364168
    "P5" = 0.5980392148364169
    "Q5" = 18.52932742631111
    "R5" = 166.7639468368
    "S5" = 0.1721768842733599
    "T5" = 0.1721768842733599
    "G6" = 0.4963216666666666
    "I6" = 0.1064050028249084
    "J6" = 0.1064050028249084
    "M6" = 4.901461666666667
    "N6" = 14.704385
    "O6" = 0.2124427850531459
    "P6" = 0.2124427850531459
    "Q6" = 2.432701623502778
    "R6" = 21.894314611525
    "S6" = 0.0226049751437114
    "T6" = 0.0226049751437114
    "G7" = 0.4963216666666666
    "I7" = 0.1064050028249084
    "J7" = 0.1064050028249084
    "O7" = 0.1372144215401173
    "P7" = 0.1372144215401173
    "Q7" = 1.571254801452222
    "S7" = 0.01460030091159435
    "T7" = 0.01460030091159435
    "G8" = 0.4963216666666666
    "I8" = 0.1064050028249084
    "J8" = 0.1064050028249084
    "M8" = 1.206743666666667
    "N8" = 3.620231
    "O8" = 0.05230357857032003
    "P8" = 0.05230357857032004
    "Q8" = 0.5989330278794444
    "R8" = 5.390397250914999
    "S8" = 0.005565362425527721
    "T8" = 0.005565362425527722
    "G9" = 0.4963216666666666
    "I9" = 0.1064050028249084
    "J9" = 0.1064050028249084
    "M9" = 13.79790933333333
    "N9" = 41.393728
    "O9" = 0.5980392148364168
    "P9" = 0.5980392148364169
    "Q9" = 6.848201356835554
    "R9" = 61.63381221151999
    "S9" = 0.06363436434407493
    "T9" = 0.06363436434407493
    "E10" = 1
    "F10" = 0.3333333333333333
    "G10" = 0.1501973333333333
    "H10" = 0.450592
    "I10" = 0.03220038283833477
    "J10" = 0.03220038283833477
    "M10" = 4.901461666666667
    "N10" = 14.704385
    "O10" = 0.2124427850531459
    "P10" = 0.2124427850531459
    "Q10" = 0.7361864717688889
    "R10" = 6.62567824592
    "S10" = 0.00684073900995336
    "T10" = 0.006840739009953361
    "E11" = 1
    "F11" = 0.3333333333333333
    "G11" = 0.1501973333333333
    "H11" = 0.450592
    "I11" = 0.03220038283833477
    "J11" = 0.03220038283833477
    "O11" = 0.1372144215401173
    "P11" = 0.1372144215401173
    "Q11" = 0.4754946177351111
    "R11" = 4.279451559616001
    "S11" = 0.004418356904532426
    "T11" = 0.004418356904532426
    "E12" = 1
    "F12" = 0.3333333333333333
    "G12" = 0.1501973333333333
    "H12" = 0.450592
    "I12" = 0.03220038283833477
    "J12" = 0.03220038283833477
    "M12" = 1.206743666666667
    "N12" = 3.620231
    "O12" = 0.05230357857032003
    "P12" = 0.05230357857032004
    "Q12" = 0.1812496807502222
    "R12" = 1.631247126752
    "S12" = 0.001684195253779227
    "T12" = 0.001684195253779228
    "E13" = 1
    "F13" = 0.3333333333333333
    "G13" = 0.1501973333333333
    "H13" = 0.450592
    "I13" = 0.03220038283833477
    "J13" = 0.03220038283833477
    "M13" = 13.79790933333333
    "N13" = 41.393728
    "O13" = 0.5980392148364168
    "P13" = 0.5980392148364169
    "Q13" = 2.072409187441777
    "R13" = 18.651682686976
    "S13" = 0.01925709167006975
    "T13" = 0.01925709167006976
    "G14" = 2.675030666666667
    "H14" = 8.025092000000001
    "I14" = 0.573492282847582
    "J14" = 0.573492282847582
    "M14" = 4.901461666666667
    "N14" = 14.704385
    "O14" = 0.2124427850531459
    "P14" = 0.2124427850531459
    "Q14" = 13.11156026982445
    "R14" = 118.00404242842
    "S14" = 0.1218342977746268
    "T14" = 0.1218342977746268
    "G15" = 2.675030666666667
    "H15" = 8.025092000000001
    "I15" = 0.573492282847582
    "J15" = 0.573492282847582
    "O15" = 0.1372144215401173
    "P15" = 0.1372144215401173
    "Q15" = 8.468610301179556
    "R15" = 76.21749271061601
    "S15" = 0.07869141184865229
    "T15" = 0.07869141184865229
    "G16" = 2.675030666666667
    "H16" = 8.025092000000001
    "I16" = 0.573492282847582
    "J16" = 0.573492282847582
    "M16" = 1.206743666666667
    "N16" = 3.620231
    "O16" = 0.05230357857032003
    "P16" = 0.05230357857032004
    "Q16" = 3.228076315139111
    "R16" = 29.052686836252
    "S16" = 0.02999569867539071
    "T16" = 0.02999569867539071
    "G17" = 2.675030666666667
    "H17" = 8.025092000000001
    "I17" = 0.573492282847582
    "J17" = 0.573492282847582
    "M17" = 13.79790933333333
    "N17" = 41.393728
    "O17" = 0.5980392148364168
    "P17" = 0.5980392148364169
    "Q17" = 36.90983060255289
    "R17" = 332.188475422976
    "S17" = 0.3429708745489122
    "T17" = 0.3429708745489122
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
